$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# add factor for state purchases: update total/impact/state_local values for
# rows 87-90 (new "state purchases" factor changes total & state_local together
# with impact), and the impact-only values for rows 91-93.

$ws.Range("E87").Value = 6.807648377777404
$ws.Range("F87").Value = 5.369168832839658
$ws.Range("H87").Value = 0.1642490337702944

$ws.Range("E88").Value = -3.020017250891981
$ws.Range("F88").Value = 1.041720123925258
$ws.Range("H88").Value = 0.2541332789429298

$ws.Range("E89").Value = -2.956577902985297
$ws.Range("F89").Value = -0.6191846308248271
$ws.Range("H89").Value = -0.1580530711298772

$ws.Range("E90").Value = -5.994976992591383
$ws.Range("F90").Value = -1.290980942172814
$ws.Range("H90").Value = 0.03569302105791428

$ws.Range("F91").Value = -4.584922296837583
$ws.Range("F92").Value = -4.323962679785499
$ws.Range("F93").Value = -3.716567364032327
